$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename odds columns (Odd_H/Odd_D/Odd_A -> H_odd/D_odd/A_odd) ---
$ws.Cells.Item(1,9).Value = 'H_odd'
$ws.Cells.Item(1,10).Value = 'D_odd'
$ws.Cells.Item(1,11).Value = 'A_odd'

# --- Extend the table: rows 13-18 are brand new, so give column A (No.) the same
# bold/bordered/centered style already used by A2:A12 before filling in values. ---
$ws.Range("A2").Copy()
$ws.Range("A13:A18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 2 (No. 1): Lausanne Ouchy vs Lausanne
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 'pzakEcuS'
$ws.Cells.Item(2,3).Value = '10/02/2024'
$ws.Cells.Item(2,4).Value = '14:00'
$ws.Cells.Item(2,5).Value = 'SWITZERLAND - SUPER LEAGUE'
$ws.Cells.Item(2,6).Value = 'Lausanne Ouchy'
$ws.Cells.Item(2,7).Value = 'Lausanne'
$ws.Cells.Item(2,8).Value = 'ROUND 23'
$ws.Cells.Item(2,9).Value = 3.4
$ws.Cells.Item(2,10).Value = 3.5
$ws.Cells.Item(2,11).Value = 2.1

# Row 3 (No. 2): Zurich vs Grasshoppers
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = 'xSjH9Hmq'
$ws.Cells.Item(3,3).Value = '10/02/2024'
$ws.Cells.Item(3,4).Value = '14:00'
$ws.Cells.Item(3,5).Value = 'SWITZERLAND - SUPER LEAGUE'
$ws.Cells.Item(3,6).Value = 'Zurich'
$ws.Cells.Item(3,7).Value = 'Grasshoppers'
$ws.Cells.Item(3,8).Value = 'ROUND 23'
$ws.Cells.Item(3,9).Value = 1.8
$ws.Cells.Item(3,10).Value = 3.6
$ws.Cells.Item(3,11).Value = 4.5

# Row 4 (No. 3): Lugano vs Young Boys
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = 'GIkL8y3k'
$ws.Cells.Item(4,3).Value = '10/02/2024'
$ws.Cells.Item(4,4).Value = '16:30'
$ws.Cells.Item(4,5).Value = 'SWITZERLAND - SUPER LEAGUE'
$ws.Cells.Item(4,6).Value = 'Lugano'
$ws.Cells.Item(4,7).Value = 'Young Boys'
$ws.Cells.Item(4,8).Value = 'ROUND 23'
$ws.Cells.Item(4,9).Value = 3
$ws.Cells.Item(4,10).Value = 3.6
$ws.Cells.Item(4,11).Value = 2.3

# Row 5 (No. 4): Winterthur vs Luzern
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = 'bsuQ7eIe'
$ws.Cells.Item(5,3).Value = '11/02/2024'
$ws.Cells.Item(5,4).Value = '10:15'
$ws.Cells.Item(5,5).Value = 'SWITZERLAND - SUPER LEAGUE'
$ws.Cells.Item(5,6).Value = 'Winterthur'
$ws.Cells.Item(5,7).Value = 'Luzern'
$ws.Cells.Item(5,8).Value = 'ROUND 23'
$ws.Cells.Item(5,9).Value = 2.63
$ws.Cells.Item(5,10).Value = 3.5
$ws.Cells.Item(5,11).Value = 2.6

# Row 6 (No. 5): Basel vs St. Gallen
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = 'CpqU6FX1'
$ws.Cells.Item(6,3).Value = '11/02/2024'
$ws.Cells.Item(6,4).Value = '12:30'
$ws.Cells.Item(6,5).Value = 'SWITZERLAND - SUPER LEAGUE'
$ws.Cells.Item(6,6).Value = 'Basel'
$ws.Cells.Item(6,7).Value = 'St. Gallen'
$ws.Cells.Item(6,8).Value = 'ROUND 23'
$ws.Cells.Item(6,9).Value = 2.8
$ws.Cells.Item(6,10).Value = 3.6
$ws.Cells.Item(6,11).Value = 2.4

# Row 7 (No. 6): Yverdon vs Servette
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = 'QmXBVbA2'
$ws.Cells.Item(7,3).Value = '11/02/2024'
$ws.Cells.Item(7,4).Value = '12:30'
$ws.Cells.Item(7,5).Value = 'SWITZERLAND - SUPER LEAGUE'
$ws.Cells.Item(7,6).Value = 'Yverdon'
$ws.Cells.Item(7,7).Value = 'Servette'
$ws.Cells.Item(7,8).Value = 'ROUND 23'
$ws.Cells.Item(7,9).Value = 5
$ws.Cells.Item(7,10).Value = 4.2
$ws.Cells.Item(7,11).Value = 1.65

# Row 8 (No. 7): Always Ready vs Universitario de Vinto
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = 'YyFMxQv1'
$ws.Cells.Item(8,3).Value = '16/02/2024'
$ws.Cells.Item(8,4).Value = '21:00'
$ws.Cells.Item(8,5).Value = 'BOLIVIA - DIVISION PROFESIONAL'
$ws.Cells.Item(8,6).Value = 'Always Ready'
$ws.Cells.Item(8,7).Value = 'Universitario de Vinto'
$ws.Cells.Item(8,8).Value = 'APERTURA '
$ws.Cells.Item(8,9).ClearContents()
$ws.Cells.Item(8,10).ClearContents()
$ws.Cells.Item(8,11).ClearContents()

# Row 9 (No. 8): Sloga Doboj vs Velez Mostar
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = '6sXRvujr'
$ws.Cells.Item(9,3).Value = '16/02/2024'
$ws.Cells.Item(9,4).Value = '14:00'
$ws.Cells.Item(9,5).Value = 'BOSNIA AND HERZEGOVINA - PREMIJER LIGA BIH'
$ws.Cells.Item(9,6).Value = 'Sloga Doboj'
$ws.Cells.Item(9,7).Value = 'Velez Mostar'
$ws.Cells.Item(9,8).Value = 'ROUND 19'
$ws.Cells.Item(9,9).ClearContents()
$ws.Cells.Item(9,10).ClearContents()
$ws.Cells.Item(9,11).ClearContents()

# Row 10 (No. 9): Lausanne vs Yverdon
$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = '4vWFUIP8'
$ws.Cells.Item(10,3).Value = '17/02/2024'
$ws.Cells.Item(10,4).Value = '14:00'
$ws.Cells.Item(10,5).Value = 'SWITZERLAND - SUPER LEAGUE'
$ws.Cells.Item(10,6).Value = 'Lausanne'
$ws.Cells.Item(10,7).Value = 'Yverdon'
$ws.Cells.Item(10,8).Value = 'ROUND 24'
$ws.Cells.Item(10,9).ClearContents()
$ws.Cells.Item(10,10).ClearContents()
$ws.Cells.Item(10,11).ClearContents()

# Row 11 (No. 10): Zeljeznicar vs Zrinjski
$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = 'zatdp1SK'
$ws.Cells.Item(11,3).Value = '17/02/2024'
$ws.Cells.Item(11,4).Value = '12:00'
$ws.Cells.Item(11,5).Value = 'BOSNIA AND HERZEGOVINA - PREMIJER LIGA BIH'
$ws.Cells.Item(11,6).Value = 'Zeljeznicar'
$ws.Cells.Item(11,7).Value = 'Zrinjski'
$ws.Cells.Item(11,8).Value = 'ROUND 19'
$ws.Cells.Item(11,9).ClearContents()
$ws.Cells.Item(11,10).ClearContents()
$ws.Cells.Item(11,11).ClearContents()

# Row 12 (No. 11): Posusje vs Zvijezda 09
$ws.Cells.Item(12,1).Value = 11
$ws.Cells.Item(12,2).Value = '29QzwJKf'
$ws.Cells.Item(12,3).Value = '17/02/2024'
$ws.Cells.Item(12,4).Value = '09:00'
$ws.Cells.Item(12,5).Value = 'BOSNIA AND HERZEGOVINA - PREMIJER LIGA BIH'
$ws.Cells.Item(12,6).Value = 'Posusje'
$ws.Cells.Item(12,7).Value = 'Zvijezda 09'
$ws.Cells.Item(12,8).Value = 'ROUND 19'
$ws.Cells.Item(12,9).ClearContents()
$ws.Cells.Item(12,10).ClearContents()
$ws.Cells.Item(12,11).ClearContents()

# Row 13 (No. 12): SA Bulo Bulo vs The Strongest
$ws.Cells.Item(13,1).Value = 12
$ws.Cells.Item(13,2).Value = 'nLCAuSfr'
$ws.Cells.Item(13,3).Value = '17/02/2024'
$ws.Cells.Item(13,4).Value = '16:00'
$ws.Cells.Item(13,5).Value = 'BOLIVIA - DIVISION PROFESIONAL'
$ws.Cells.Item(13,6).Value = 'SA Bulo Bulo'
$ws.Cells.Item(13,7).Value = 'The Strongest'
$ws.Cells.Item(13,8).Value = 'APERTURA '
$ws.Cells.Item(13,9).ClearContents()
$ws.Cells.Item(13,10).ClearContents()
$ws.Cells.Item(13,11).ClearContents()

# Row 14 (No. 13): Santa Cruz vs Tomayapo
$ws.Cells.Item(14,1).Value = 13
$ws.Cells.Item(14,2).Value = 'tvzrp8nR'
$ws.Cells.Item(14,3).Value = '17/02/2024'
$ws.Cells.Item(14,4).Value = '21:00'
$ws.Cells.Item(14,5).Value = 'BOLIVIA - DIVISION PROFESIONAL'
$ws.Cells.Item(14,6).Value = 'Santa Cruz'
$ws.Cells.Item(14,7).Value = 'Tomayapo'
$ws.Cells.Item(14,8).Value = 'APERTURA '
$ws.Cells.Item(14,9).ClearContents()
$ws.Cells.Item(14,10).ClearContents()
$ws.Cells.Item(14,11).ClearContents()

# Row 15 (No. 14): Bolivar vs GV San Jose
$ws.Cells.Item(15,1).Value = 14
$ws.Cells.Item(15,2).Value = 'r7dDamnE'
$ws.Cells.Item(15,3).Value = '17/02/2024'
$ws.Cells.Item(15,4).Value = '18:30'
$ws.Cells.Item(15,5).Value = 'BOLIVIA - DIVISION PROFESIONAL'
$ws.Cells.Item(15,6).Value = 'Bolivar'
$ws.Cells.Item(15,7).Value = 'GV San Jose'
$ws.Cells.Item(15,8).Value = 'APERTURA '
$ws.Cells.Item(15,9).ClearContents()
$ws.Cells.Item(15,10).ClearContents()
$ws.Cells.Item(15,11).ClearContents()

# Row 16 (No. 15): St. Gallen vs Winterthur
$ws.Cells.Item(16,1).Value = 15
$ws.Cells.Item(16,2).Value = 'pjmDTxuF'
$ws.Cells.Item(16,3).Value = '17/02/2024'
$ws.Cells.Item(16,4).Value = '14:00'
$ws.Cells.Item(16,5).Value = 'SWITZERLAND - SUPER LEAGUE'
$ws.Cells.Item(16,6).Value = 'St. Gallen'
$ws.Cells.Item(16,7).Value = 'Winterthur'
$ws.Cells.Item(16,8).Value = 'ROUND 24'
$ws.Cells.Item(16,9).ClearContents()
$ws.Cells.Item(16,10).ClearContents()
$ws.Cells.Item(16,11).ClearContents()

# Row 17 (No. 16): Igman K. vs GOSK Gabela
$ws.Cells.Item(17,1).Value = 16
$ws.Cells.Item(17,2).Value = 'l2p0qLsR'
$ws.Cells.Item(17,3).Value = '17/02/2024'
$ws.Cells.Item(17,4).Value = '09:00'
$ws.Cells.Item(17,5).Value = 'BOSNIA AND HERZEGOVINA - PREMIJER LIGA BIH'
$ws.Cells.Item(17,6).Value = 'Igman K.'
$ws.Cells.Item(17,7).Value = 'GOSK Gabela'
$ws.Cells.Item(17,8).Value = 'ROUND 19'
$ws.Cells.Item(17,9).ClearContents()
$ws.Cells.Item(17,10).ClearContents()
$ws.Cells.Item(17,11).ClearContents()

# Row 18 (No. 17): Grasshoppers vs Basel
$ws.Cells.Item(18,1).Value = 17
$ws.Cells.Item(18,2).Value = 'OtnHSdfL'
$ws.Cells.Item(18,3).Value = '17/02/2024'
$ws.Cells.Item(18,4).Value = '16:30'
$ws.Cells.Item(18,5).Value = 'SWITZERLAND - SUPER LEAGUE'
$ws.Cells.Item(18,6).Value = 'Grasshoppers'
$ws.Cells.Item(18,7).Value = 'Basel'
$ws.Cells.Item(18,8).Value = 'ROUND 24'
$ws.Cells.Item(18,9).ClearContents()
$ws.Cells.Item(18,10).ClearContents()
$ws.Cells.Item(18,11).ClearContents()

